# Updates the cryptos worksheet with the latest scraped price/volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value reads as a plain number (e.g. "300.22") need the
# Text number format pre-applied, otherwise Excel COM auto-converts the
# assignment into a numeric cell instead of keeping it as text like the source data.
$textForceRefs = @('D5', 'D6', 'D7', 'D9', 'D10', 'D11', 'D12', 'D16', 'D17', 'D19', 'D21', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27', 'D29', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D36', 'D37', 'D39', 'D40', 'D41', 'D42', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D51')
foreach ($ref in $textForceRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range('D2').Value = '45.597.19'
$ws.Range('E2').Value = '  -2.04%  '
$ws.Range('D3').Value = '2.415.34'
$ws.Range('E3').Value = '  +5.34%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '300.22'
$ws.Range('E5').Value = '  -1.39%  '
$ws.Range('D6').Value = '97.55'
$ws.Range('E6').Value = '  -3.41%  '
$ws.Range('D7').Value = '0.563'
$ws.Range('E7').Value = '  -0.57%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('D9').Value = '0.509'
$ws.Range('E9').Value = '  -1.93%  '
$ws.Range('D10').Value = '34.56'
$ws.Range('E10').Value = '  -5.78%  '
$ws.Range('D11').Value = '0.0791'
$ws.Range('E11').Value = '  +0.07%  '
$ws.Range('D12').Value = '7.16'
$ws.Range('E12').Value = '  -2.65%  '
$ws.Range('E13').Value = '  +0.59%  '
$ws.Range('D14').Value = '2.768.24'
$ws.Range('E14').Value = '  +4.83%  '
$ws.Range('D15').Value = '2.407.70'
$ws.Range('E15').Value = '  +5.19%  '
$ws.Range('D16').Value = '14.16'
$ws.Range('E16').Value = '  +1.88%  '
$ws.Range('D17').Value = '0.832'
$ws.Range('E17').Value = '  +2.56%  '
$ws.Range('D18').Value = '45.574.88'
$ws.Range('E18').Value = '  -2.09%  '
$ws.Range('D19').Value = '12.95'
$ws.Range('E19').Value = '  -0.97%  '
$ws.Range('D20').Value = '0.0₃0952'
$ws.Range('E20').Value = '  +1.55%  '
$ws.Range('D21').Value = '6.13'
$ws.Range('E21').Value = '  +1.64%  '
$ws.Range('D22').Value = '67.37'
$ws.Range('E22').Value = '  +1.62%  '
$ws.Range('D23').Value = '243.51'
$ws.Range('E23').Value = '  -1.65%  '
$ws.Range('D24').Value = '2.79'
$ws.Range('E24').Value = '  -4.17%  '
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('D26').Value = '1.93'
$ws.Range('E26').Value = '  +0.38%  '
$ws.Range('D27').Value = '39.23'
$ws.Range('E27').Value = '  -9.25%  '
$ws.Range('E28').Value = '  -1.99%  '
$ws.Range('D29').Value = '9.77'
$ws.Range('E29').Value = '  -0.63%  '
$ws.Range('D30').Value = '3.85'
$ws.Range('E30').Value = '  +19.21%  '
$ws.Range('D31').Value = '21.40'
$ws.Range('E31').Value = '  +7.28%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').Value = '5.57'
$ws.Range('E32').Value = '  -1.27%  '
$ws.Range('B33').Value = 'WEMIXToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D33').Value = '2.74'
$ws.Range('E33').Value = '  -2.27%  '
$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D34').Value = '147.71'
$ws.Range('E34').Value = '  +0.60%  '
$ws.Range('D35').Value = '0.0776'
$ws.Range('E35').Value = '  -2.38%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = '1.97'
$ws.Range('E36').Value = '  +11.19%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').Value = '0.112'
$ws.Range('E37').Value = '  -2.62%  '
$ws.Range('E38').Value = '  -1.02%  '
$ws.Range('D39').Value = '15.50'
$ws.Range('E39').Value = '  -2.92%  '
$ws.Range('D40').Value = '3.88'
$ws.Range('E40').Value = '  -4.39%  '
$ws.Range('D41').Value = '0.0300'
$ws.Range('E41').Value = '  -1.05%  '
$ws.Range('D42').Value = '3.27'
$ws.Range('E42').Value = '  -3.18%  '
$ws.Range('D43').Value = '1.955.66'
$ws.Range('E43').Value = '  +6.64%  '
$ws.Range('D44').Value = '0.998'
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').Value = '92.53'
$ws.Range('E45').Value = '  +5.91%  '
$ws.Range('D46').Value = '1.78'
$ws.Range('E46').Value = '  -10.17%  '
$ws.Range('D47').Value = '8.65'
$ws.Range('E47').Value = '  +10.30%  '
$ws.Range('D48').Value = '99.75'
$ws.Range('E48').Value = '  +4.26%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = '0.185'
$ws.Range('E49').Value = '  -4.90%  '
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '2.637.50'
$ws.Range('E50').Value = '  +4.78%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '14.50'
$ws.Range('E51').Value = '  +7.37%  '

# Drop back to the default style now that the text is safely stored, so the
# cells don't carry a stray Text-format style like a real typed apostrophe would.
foreach ($ref in $textForceRefs) {
    $ws.Range($ref).Style = "Normal"
}
